$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet - row 3 is the b.md entry. Its zh-cn / de-de status moves
# from "Handed back: in sync with en-US" to "Ready for handoff", and the
# "Latest HO Xliff Generate Date" bumps to the new handoff timestamp.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-17 15:42:24"

# ---------------------------------------------------------------------------
# "zh-cn" sheet - row 3 (b.md) gets a fresh handoff: new status, new handoff
# file name / datetime, and a new "stale handback" error detail message.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-10-17 15:42:03"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2f6d5ac46497618ff02a23baca4756261bd7666/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/717ea06fbb6b10fab1c10ad77460e2539d8f9684/e2e/b.md."
# Column P (Error Detail) widens to fit the long message.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# "de-de" sheet - same kind of update for row 3 (b.md).
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-10-17 15:42:24"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2f6d5ac46497618ff02a23baca4756261bd7666/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/717ea06fbb6b10fab1c10ad77460e2539d8f9684/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
